# comment heartbeat + unit test patient
#
# Adds a new "unit test patient" minute-average reading to the
# "Vital Signs Average per Minute" sheet, and records four abnormal-event
# rows (RespRate, ECG x2, Temperature) on the "Abnormal Events" sheet,
# replacing the previous single ECG "heartbeat" comment event.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Vital Signs Average per Minute
$ws2 = $wb.Worksheets.Item(2)   # Abnormal Events

# --- Sheet 1: Vital Signs Average per Minute -------------------------------
# New row 2: a single averaged-minute vitals reading for the test patient.
$ws1.Range("A2").Value = "2026-01-14 01:42"
$ws1.Range("B2").Value = 75.0
$ws1.Range("C2").Value = 13.0
$ws1.Range("D2").Value = 36.0
$ws1.Range("E2").Value = "116/83"

# --- Sheet 2: Abnormal Events ------------------------------------------------
# Row 2: RespRate abnormal event
$ws2.Range("A2").Value = "2026-01-14 01:43:05"
$ws2.Range("B2").Value = "2026-01-14 01:43:05"
$ws2.Range("C2").Value = "RespRate"
$ws2.Range("D2").Value = "11.0 - 11.0"
$ws2.Range("E2").Value = "AMBER"

# Row 3: ECG abnormal event
$ws2.Range("A3").Value = "2026-01-14 01:43:06"
$ws2.Range("B3").Value = "2026-01-14 01:43:06"
$ws2.Range("C3").Value = "ECG"
$ws2.Range("D3").Value = "0.7 - 0.7"
$ws2.Range("E3").Value = "AMBER"

# Row 4: second ECG abnormal event
$ws2.Range("A4").Value = "2026-01-14 01:43:09"
$ws2.Range("B4").Value = "2026-01-14 01:43:10"
$ws2.Range("C4").Value = "ECG"
$ws2.Range("D4").Value = "-0.8 - 0.8"
$ws2.Range("E4").Value = "AMBER"

# Row 5: Temperature abnormal event
$ws2.Range("A5").Value = "2026-01-14 01:43:12"
$ws2.Range("B5").Value = "2026-01-14 01:43:13"
$ws2.Range("C5").Value = "Temperature"
$ws2.Range("D5").Value = "35.8 - 35.9"
$ws2.Range("E5").Value = "AMBER"
